$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Amazon"
$ws.Range("B3").Value = "Ebay"
$ws.Range("B4").Value = "Wish"

$ws.Range("B5").Select()
